$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '66.447.85'
Set-TextValue 'E2' '  +0.07%  '
Set-TextValue 'D3' '3.298.42'
Set-TextValue 'E3' '  -0.67%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '586.77'
Set-TextValue 'E5' '  +1.74%  '
Set-TextValue 'D6' '180.17'
Set-TextValue 'E6' '  -0.50%  '
Set-TextValue 'D7' '0.642'
Set-TextValue 'E7' '  +1.62%  '
Set-TextValue 'E8' '  +0.05%  '
Set-TextValue 'D9' '3.295.35'
Set-TextValue 'E9' '  -0.55%  '
Set-TextValue 'E10' '  -1.50%  '
Set-TextValue 'D11' '6.85'
Set-TextValue 'E11' '  +2.21%  '
Set-TextValue 'E12' '  -0.60%  '
Set-TextValue 'D13' '3.876.64'
Set-TextValue 'E13' '  -0.41%  '
Set-TextValue 'E14' '  -2.33%  '
Set-TextValue 'D15' '66.422.66'
Set-TextValue 'E15' '  -0.05%  '
Set-TextValue 'D16' '26.61'
Set-TextValue 'E16' '  -0.61%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '3.306.21'
Set-TextValue 'E17' '  +0.75%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D18' '0.0000163'
Set-TextValue 'E18' '  -1.44%  '
Set-TextValue 'D19' '427.42'
Set-TextValue 'E19' '  -3.08%  '
Set-TextValue 'D20' '5.44'
Set-TextValue 'E20' '  -3.75%  '
Set-TextValue 'D21' '13.02'
Set-TextValue 'E21' '  -3.79%  '
Set-TextValue 'D22' '7.32'
Set-TextValue 'E22' '  -2.95%  '
Set-TextValue 'D23' '71.59'
Set-TextValue 'E23' '  -2.40%  '
Set-TextValue 'D24' '0.999'
Set-TextValue 'E24' '  -0.10%  '
Set-TextValue 'D25' '5.73'
Set-TextValue 'E25' '  +0.66%  '
Set-TextValue 'D26' '3.455.87'
Set-TextValue 'D27' '0.515'
Set-TextValue 'E27' '  -1.11%  '
Set-TextValue 'D28' '0.206'
Set-TextValue 'E28' '  +4.89%  '
Set-TextValue 'D29' '0.0000114'
Set-TextValue 'E29' '  -1.36%  '
Set-TextValue 'D30' '9.11'
Set-TextValue 'E30' '  +0.37%  '
Set-TextValue 'E31' '  +0.13%  '
Set-TextValue 'D32' '1.91'
Set-TextValue 'E32' '  -1.90%  '
Set-TextValue 'D33' '22.33'
Set-TextValue 'E33' '  -1.78%  '
Set-TextValue 'D35' '5.16'
Set-TextValue 'E35' '  -1.06%  '
Set-TextValue 'D36' '6.57'
Set-TextValue 'E36' '  -2.83%  '
Set-TextValue 'D37' '1.18'
Set-TextValue 'E37' '  -2.57%  '
Set-TextValue 'D38' '158.87'
Set-TextValue 'E38' '  -0.20%  '
Set-TextValue 'E39' '  -3.21%  '
Set-TextValue 'D40' '2.873.47'
Set-TextValue 'E40' '  +1.53%  '
Set-TextValue 'E41' '  -1.28%  '
Set-TextValue 'D42' '26.26'
Set-TextValue 'E42' '  -3.88%  '
Set-TextValue 'D43' '4.32'
Set-TextValue 'E43' '  -2.64%  '
Set-TextValue 'D44' '0.751'
Set-TextValue 'E44' '  -4.72%  '
Set-TextValue 'D45' '39.68'
Set-TextValue 'E45' '  -2.33%  '
Set-TextValue 'D46' '0.0656'
Set-TextValue 'E46' '  -1.63%  '
Set-TextValue 'D47' '5.91'
Set-TextValue 'E47' '  -4.23%  '
Set-TextValue 'D48' '2.30'
Set-TextValue 'E48' '  -1.45%  '
Set-TextValue 'D49' '22.92'
Set-TextValue 'E49' '  -4.63%  '
Set-TextValue 'D50' '310.84'
Set-TextValue 'E50' '  -5.19%  '
Set-TextValue 'E51' '  -0.72%  '
